# fixed #741 FlixelRL-741 STATSの説明ポップアップ表示
#
# Insert a new "stats" worksheet right before the "deathtype" sheet.
# It reuses the exact row/column layout & styling of "deathtype" (same
# table shape: id/msg header, int/str type row, 6 numbered entries), so
# the cleanest way to reproduce that formatting through the Excel object
# model is to duplicate "deathtype" and then overwrite the label texts.

$wb = $excel.ActiveWorkbook

$deathtype = $wb.Worksheets.Item("deathtype")

# Copy "deathtype" to a position right before itself, then rename the
# copy to "stats". Resulting order: ... title, stats, deathtype,
# staffroll, statistics  (everything after "title" keeps its original
# content, just shifted one slot to the right).
$deathtype.Copy($deathtype)
$stats = $wb.Worksheets.Item("deathtype (2)")
$stats.Name = "stats"

# Overwrite the 6 data rows with the new STATS popup description text.
$stats.Range("B3").Value = "プレイ情報を見ます"
$stats.Range("B4").Value = "実績データを閲覧します"
$stats.Range("B5").Value = "ゲームプレイ履歴を見ます"
$stats.Range("B6").Value = "倒した敵の情報を見ます"
$stats.Range("B7").Value = "獲得したアイテムを見ます"
$stats.Range("B8").Value = "タイトル画面に戻ります"
